# Auto-generated COM-interop script recreating the product_template.xlsx edit.
# Adds 5 new product rows (24-28) to the single worksheet, with matching
# shared-string / row-height / view-pane updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights for the new rows ---
$ws.Rows.Item(24).RowHeight = 375
$ws.Rows.Item(25).RowHeight = 405
$ws.Rows.Item(26).RowHeight = 375
$ws.Rows.Item(27).RowHeight = 345
$ws.Rows.Item(28).RowHeight = 409.5

# --- Numeric cells (price, quantity) ---
$ws.Range("D24").Value = 250
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 300
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 250
$ws.Range("E26").Value = 0
$ws.Range("D27").Value = 300
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 140
$ws.Range("E28").Value = 0

# --- Text cells, written in the exact order the source file introduced
#     them so the rebuilt shared-string table lines up with the target ---
$v24A = @'
AA SKIN BOOST HA+ Day/night cream 50ml
'@
$ws.Range("A24").Value = $v24A

$v24B = @'
 Face Moisturizing
'@
$ws.Range("B24").Value = $v24B

$v24C = @'
Description:
AA SKIN BOOST Day and night HYDRO-LOCK creamrich in hyaluronic acid, enhanced with lipid protection system and vitamin E, is recommended for intense daily moisturising and nourishing. You can add a few drops of the AA Skin Boost Concentrate with hyaluronic acid to strengthen the effect.
Effect
Deeply hydrated, nourished and soft to the touch skin, restored smoothness and firmness, healthy and radiant look.
'@
$ws.Range("C24").Value = $v24C

$v24F = @'
ايه ايه اسكن بوست اتش ايه بلس داي / نايت
   كريم الترطيب المركز والتغذية للبشرة بالهيالورونيك اسيد
(كريم الوجه مناسب لكل أنواع البشرة ومختبر من قبل أطباء الجلد والحساسية) 

'@
$ws.Range("F24").Value = $v24F

$v24G = @'
الوصف:                                                                  
يعمل كريم ايه ايه اسكن بوست اتش ايه بلس  بنظام الهايدرولوك المدعم بالهيالورونيك اسيد وفيتامين إي علي تغذية وترطيب البشرة بطريقه مكثفه مما يساعد البشرة علي استعادة النعومة والثبات والتمتع بمظهر صحي ومشرق

'@
$ws.Range("G24").Value = $v24G

$v25F = @'
ايه ايه اسكن بوست هيالورونيك أسيد بلس كومبليكس 15% سيروم الترطيب المركز و الحماية من العوامل الخارجية
(سيروم الوجه مناسب لكل أنواع البشرة ومختبر من قبل أطباء الجلد والحساسية)

'@
$ws.Range("F25").Value = $v25F

$v25A = @'
AA SKIN BOOST HA+ complex 15% Concentrate 30ml
'@
$ws.Range("A25").Value = $v25A

$v25C = @'
Description:
Concentrate with hyaluronic acid and algae extract is an active moisturizing and detoxifying treatment, strengthening the effects od daily care. Effectively smooths and revitalizes fine wrinkled.
Effect
Deeply hydrated skin, improoved elasticity, smooth wrinkle-free skin, radiant, revitalized complexion.

'@
$ws.Range("C25").Value = $v25C

$v25G = @'
الوصف:
يعمل سيروم الهيالورينيك أسيد المركز وخلاصه طحالب ألجانيا روبنز علي الترطيب العميق للبشرة مما يجعلها أكثر نعومه ومرونة لتقاوم ظهور التجاعيد وتكون البشرة مليئة بالحيوية .

'@
$ws.Range("G25").Value = $v25G

$v26A = @'
AA SKIN BOOST C+ Day/Night Cream 50ml
'@
$ws.Range("A26").Value = $v26A

$v26C = @'
Description:
AA Skin Boost Day and night C-forte cream with collagen and hyaluronic acid isrecommended for daily antioxidant care providing skin lightening and tone evening effect. It will make your skin look healthier, radiant and more hydrated. You can add a few drops of the AA Skin Boost Concentrate with vitamin C to strenghten the effects of the cream.
Effect
Radiant, full of energy skin, healthy, fresh and rested appearance, skin protected from harmful external factors, restored optimal hydration level.
'@
$ws.Range("C26").Value = $v26C

$v26F = @'
ايه ايه اسكن بوست سي بلس داي / نايت كريم بفيتامين سي لبشره مشرقه ومليئة بالطاقة
(كريم الوجه مناسب لكل أنواع البشرة - مختبر من قبل أطباء الجلد والحساسية)

'@
$ws.Range("F26").Value = $v26F

$v26G = @'
الوصف:
يعمل كريم ايه ايه اسكن بوست سي بلس بنظام سي فورت المدعم بالكولاجين والهيالورونيك أسيد  علي تفتيح و إعطاء البشرة مظهر صحي و منعش مع حمايتها من العوامل الخارجية الضارة واستعادة مستوى الترطيب الأمثل لبشرة مشرقة مليئة بالحيوية تحتوي تركيبة الكريم علي اربعة اشكال من فيتامين سي للقوه الكاملة.  

'@
$ws.Range("G26").Value = $v26G

$v27F = @'
ايه ايه اسكن بوست سي بلس  8% سيروم فيتامين سي مع مستخلص الأسيلورا لنضارة وتوحيد لون البشرة
(سيروم الوجه مناسب لكل أنواع البشرة - مختبر من قبل أطباء الجلد والحساسية)

'@
$ws.Range("F27").Value = $v27F

$v27A = @'
AA SKIN BOOST C+ 8% Concentrate           Vitamin C  + Acelora extract 30ml
'@
$ws.Range("A27").Value = $v27A

$v27C = @'
Description:
AA Skin Boost Concentrate with vitamin C and acerola is an intensive antioxidant and brightening treatment with skin brightening and tone evening effect, preventing premature skin ageing and strengthening the effects of daily care.
Bright and radiant complexion, even skin tone, perfect firmness and smopthness, optimal hydration level, health, rested look.
'@
$ws.Range("C27").Value = $v27C

$v27G = @'
الوصف:
يعمل السيروم  بشكل مكثف كمضاد للأكسدة على نضارة وتفتيح البشرة مع  توحيد لونها لتبدو اكثر حيوية مع الحفاظ علي مستوي الترطيب المثالي مما يجعلها أكثر مرونة وتقاوم ظهور الشيخوخة المبكرة

'@
$ws.Range("G27").Value = $v27G

$v28F = @'
لونج فورهير شامبو التعزيز ضد تساقط الشعر
(لكل أنواع الشعر ومناسب للرجال والنساء)

'@
$ws.Range("F28").Value = $v28F

$v28G = @'
الوصف:
يعمل الشامبو علي تقوية الشعر الضعيف والهش ويمنع فقدانه بسبب التغيرات الهرمونية والجينية ,التعب, الإجهاد ,النظام الغذائي غير المتوازن أو العوامل الخارجية.
يحتوي الشامبو علي تركيبة التيترا اوكسيدليوم     
وهو مركب مبتكر يغذي بصيلات الشعر ويمنع تساقطها, بالإضافة إلي وجود الاّلانتوين الذي يلطف التهيج ويهدئ فروة الرأس الحساسة.
يحسن شامبو لونج فورلاشز حالة الشعر بشكل واضح ويقلل من التقصف ويعطي الشعر حيوية ولمعانا

'@
$ws.Range("G28").Value = $v28G

$v28C = @'
Description:
Shampoo strengthens weak and brittle hair as well as prevents its loss caused by hormonal and genetic changes, fatigue, stress, unbalanced diet or external factors. It formulation contains Tetraxidylum – an innovative ingredients complex which nourishes hair follicles and prevents its falling out, while stimulating microcirculation. Allantoin soothes irritations and calms sensitive scalp. Shampoo visibly improves hair condition, reduces brittleness and gives it vitality and shine. 

'@
$ws.Range("C28").Value = $v28C

$v28B = @'
HAIR SHAMPOO
'@
$ws.Range("B28").Value = $v28B

$v28A = @'
LONG4HAIR by Oceanic, HAIR, Anti-hair Loss Strengthening Shampoo, 200 ml
'@
$ws.Range("A28").Value = $v28A

# --- Cells re-using the existing "Face Moisturizing" shared string ---
$ws.Range("B25").Value = "Face Moisturizing"
$ws.Range("B26").Value = "Face Moisturizing"
$ws.Range("B27").Value = "Face Moisturizing"

# --- Update frozen-pane / selection to match the edit ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("A28").Select()

Write-Host "Added rows 24-28"
